# The deck's single custom design ("Integral") is switched back to the
# default Office Theme color palette (Design > Variants > color swap),
# while the font scheme, line/fill/effect formats stay the same (the
# two themes already shared an identical font/format scheme).
#
# RGB() isn't available in this host, so the target sRGB colors are
# passed as pre-computed long (BGR-packed) integers:
#   dk1      000000 ->        0
#   lt1      FFFFFF -> 16777215
#   dk2      44546A ->  6968388
#   lt2      E7E6E6 -> 15132391
#   accent1  5B9BD5 -> 13998939
#   accent2  ED7D31 ->  3243501
#   accent3  A5A5A5 -> 10855845
#   accent4  FFC000 ->    49407
#   accent5  4472C4 -> 12874308
#   accent6  70AD47 ->  4697456
#   hlink    0563C1 -> 12673797
#   folHlink 954F72 ->  7491477

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

$colorScheme.Item(1).RGB = 0
$colorScheme.Item(2).RGB = 16777215
$colorScheme.Item(3).RGB = 6968388
$colorScheme.Item(4).RGB = 15132391
$colorScheme.Item(5).RGB = 13998939
$colorScheme.Item(6).RGB = 3243501
$colorScheme.Item(7).RGB = 10855845
$colorScheme.Item(8).RGB = 49407
$colorScheme.Item(9).RGB = 12874308
$colorScheme.Item(10).RGB = 4697456
$colorScheme.Item(11).RGB = 12673797
$colorScheme.Item(12).RGB = 7491477
